# ---------------------------------------------------------------------------
# HistoriasDeUsuarios.xlsx - "Add files via upload" re-edit
#
# The sheet "Pontos de História" is reworked:
#   * the old layout paired every story across two merged rows
#     (B3:B4, D3:D4, ...); it is collapsed into one row per story.
#   * a small reference table (Pontos de História / Complexidade /
#     Justificativa) is added at G2:I5.
#   * a print area covering both blocks is defined for this sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pontos de História")

# ---------------------------------------------------------------------------
# 1) Collapse the merged HU-pairs back into a single row per story.
#    Every story currently spans two rows (e.g. B3:B4) merged together;
#    unmerge everything then drop the now-empty second row of each pair.
# ---------------------------------------------------------------------------
$ws.Cells.UnMerge()

$blankRows = @(20, 18, 16, 14, 12, 10, 8, 6, 4)
foreach ($r in $blankRows) {
    $ws.Rows.Item($r).Delete()
}

# Header "Pontos de História" (E2) is now centered like the rest of the
# header row instead of left aligned.
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 2) New reference table: Pontos de História / Complexidade / Justificativa
# ---------------------------------------------------------------------------

# Headers - reuse the look of the existing shaded header row.
$ws.Range("C2").Copy()
$ws.Range("G2:H2").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("G2").Value = "Pontos de História"
$ws.Range("H2").Value = "Complexidade"
$ws.Range("I2").Value = "Justificativa"

# Body - reuse the existing bordered cell looks already used elsewhere on
# this sheet (plain left/indent look for G, wrap/indent look for H & I).
$ws.Range("B3").Copy()
$ws.Range("G3:G5").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("H3:H5").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("I3:I5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("G3").Value = "1 a 3"
$ws.Range("H3").Value = "Simples"
$ws.Range("I3").Value = "Processos que exigem a criação de formulários, inserções em banco e/ou implementação de módulos já existentes"

$ws.Range("G4").Value = "4 a 6"
$ws.Range("H4").Value = "Média"
$ws.Range("I4").Value = "Processos que exigem uma autenticação e/ou segurança maior e/ou regras de negócio mais complexas"

$ws.Range("G5").Value = "7 a 8"
$ws.Range("H5").Value = "Complexa"
$ws.Range("I5").Value = "Processos que exigem critérios estruturados e diagramação de banco de dados complexa, além de fluxos maiores do usuário e fluxos que devem seguir regras estritas, tanto as validando para evitar bypass pelo usuário"

# Give the new columns sensible widths.
$ws.Range("G1:H1").ColumnWidth = 15.71
$ws.Range("I1").ColumnWidth = 90.71
$ws.Range("F1").ColumnWidth = 1.71
$ws.Range("E1").ColumnWidth = 17.43

# ---------------------------------------------------------------------------
# 3) Top row of the collapsed table (row 3) gets its own boxed-without-
#    bottom-border look, separating it visually from the header above.
# ---------------------------------------------------------------------------
$topRow = $ws.Range("B3:E3")
$topRow.Borders.Item(7).LineStyle = 1
$topRow.Borders.Item(7).Weight = 2
$topRow.Borders.Item(10).LineStyle = 1
$topRow.Borders.Item(10).Weight = 2
$topRow.Borders.Item(8).LineStyle = 1
$topRow.Borders.Item(8).Weight = 2
$topRow.Borders.Item(9).LineStyle = 0
$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("D3:E3").HorizontalAlignment = -4108
$ws.Range("D3:E3").VerticalAlignment = -4108
$ws.Range("D3:E3").Font.Size = 10
$ws.Range("D3:E3").WrapText = $true

# ---------------------------------------------------------------------------
# 4) Sheet level bits: selection, print area, page breaks.
# ---------------------------------------------------------------------------
$ws.Range("C7").Select()

$ws.PageSetup.PrintArea = '$B$2:$E$11,$G$2:$I$5'
$ws.PageSetup.Zoom = $false
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1

$ws.HPageBreaks.Add($ws.Range("B2"))
$ws.VPageBreaks.Add($ws.Range("G1"))

# ---------------------------------------------------------------------------
# 5) Workbook level: define the print-area name explicitly (also produced by
#    the PageSetup.PrintArea assignment above, kept here for clarity/safety).
# ---------------------------------------------------------------------------
$wb.Names.Add("_xlnm.Print_Area", '=''Pontos de História''!$B$2:$E$11,''Pontos de História''!$G$2:$I$5')
